$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text changes (row 11) -- column headers get a "1dp"/"3sd" suffix
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "p1'/mV 1dp"
$ws.Range("C11").Value = "T1/mV 1dp"
$ws.Range("D11").Value = "p2'/mV 1dp"
$ws.Range("E11").Value = "T2/mV 1dp"
$ws.Range("F11").Value = "p1'-p2'/mV 1dp"
$ws.Range("G11").Value = "γ=p1'/p1'-p2' 3sd"
$ws.Range("H11").Value = "γ平均 3sd"
$ws.Range("I11").Value = "相对误差μ 3sd"

# ---------------------------------------------------------------------------
# 2. New rows 25/26 with labels + formulas (set before D23 so shared-string
#    append order mirrors the original authoring order)
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = "保留位数后(与μD对齐位数)"
$ws.Range("F25").Value = "保留位数后(最高次有效数字为1,2保留两位,反之保留一位)"

$ws.Range("D26").Formula = "=IF(F24*10^INT(-LOG(ABS(F24)))<0.3,FIXED(D24,1-INT(LOG(F24)),1),FIXED(D24,-INT(LOG(F24)),1))"
$ws.Range("F26").Formula = "=IF(F24*10^INT(-LOG(ABS(F24)))<0.3,FIXED(F24,1-INT(LOG(F24)),1),FIXED(F24,-INT(LOG(F24)),1))"

# Row 23 label change
$ws.Range("D23").Value = "γ平均"

# ---------------------------------------------------------------------------
# 3. Number formats: data columns B:E and F get a fixed "0.0" display format
# ---------------------------------------------------------------------------
$ws.Range("B12:E21").NumberFormat = "0.0_ "
$ws.Range("F12:F21").NumberFormat = "0.0_ "

# G:I (gamma, average, relative error) columns use scientific 2-decimal format
$ws.Range("G12:I21").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# 4. I14 loses its fill (becomes plain white) -- matches the author's manual tweak
# ---------------------------------------------------------------------------
$ws.Range("I14").Interior.ColorIndex = -4142

# ---------------------------------------------------------------------------
# 5. New cell style area (rows 25/26) formatting
# ---------------------------------------------------------------------------
$ws.Range("D25:F25").Interior.ColorIndex = -4142
$ws.Range("D26:F26").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 6. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item("B:C").ColumnWidth = 12.6640625
$ws.Columns.Item("D").ColumnWidth = 14.33203125
$ws.Columns.Item("E").ColumnWidth = 13.5546875
$ws.Columns.Item("F").ColumnWidth = 14.109375
$ws.Columns.Item("G").ColumnWidth = 16.6640625
$ws.Columns.Item("H").ColumnWidth = 15.5546875
$ws.Columns.Item("I").ColumnWidth = 19.21875

# ---------------------------------------------------------------------------
# 7. View state: scrolled down to row 3, selection on F24
# ---------------------------------------------------------------------------
$ws.Range("A3").Select()
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("F24").Select()
